$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.722773671150208
$ws.Range("B1").Value = 1.953958868980408
$ws.Range("C1").Value = 2.009294509887695
$ws.Range("D1").Value = 2.605636835098267
$ws.Range("E1").Value = 3.100084781646729
